# feat: add 2022-Q3 data
#
#  - Insert a new worksheet "2022-Q3" positioned before the existing "2022-Q2" sheet
#  - Populate it with the Q3 fund-holding data
#  - Update the "总计" (summary) sheet: rename the existing Q2 row to Q3 and append a
#    fresh row preserving the original Q2 summary figures
#
# NOTE: worksheet object variables in this environment re-resolve by their
# original positional index, so a variable captured for a sheet whose index
# shifts (because a new sheet got inserted before it) will start pointing at
# whatever sheet now occupies that old index. To stay safe we only use the
# "$wsQ2" handle to tell Add() where to insert, and otherwise always fetch
# worksheets freshly by name after any sheet-insertion happens.

$wb = $excel.ActiveWorkbook
$wsTotal = $wb.Worksheets.Item(1)
$wsQ2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# 1) Create the new "2022-Q3" sheet right before the "2022-Q2" sheet
# ---------------------------------------------------------------------------
$wsQ3 = $wb.Worksheets.Add($wsQ2)
$wsQ3.Name = "2022-Q3"

# Copy header formatting (bold/centered/bordered) from the "总计" sheet header
# row so the new sheet's header matches the style used for this workbook's
# freshly generated sheets.
$wsTotal.Range("B1:D1").Copy()
$wsQ3.Range("B1:H1").PasteSpecial(-4122)

# Copy the "index" column style (A2 on 总计) down column A for the 4 data rows
$wsTotal.Range("A2").Copy()
$wsQ3.Range("A2:A5").PasteSpecial(-4122)

# Header labels
$wsQ3.Cells.Item(1, 2).Value = "基金代码"
$wsQ3.Cells.Item(1, 3).Value = "基金名称"
$wsQ3.Cells.Item(1, 4).Value = "基金规模"
$wsQ3.Cells.Item(1, 5).Value = "股票总仓位"
$wsQ3.Cells.Item(1, 6).Value = "仓位占比"
$wsQ3.Cells.Item(1, 7).Value = "持有市值(亿元)"
$wsQ3.Cells.Item(1, 8).Value = "仓位排名"

# Force columns B:G (fund code / name / figures that must stay textual, e.g.
# to keep leading zeros in fund codes) to be stored as text before typing
# their values.
$wsQ3.Range("B2:G5").NumberFormat = "@"

# Row 2
$wsQ3.Cells.Item(2, 1).Value = 0
$wsQ3.Cells.Item(2, 2).Value = "002068"
$wsQ3.Cells.Item(2, 3).Value = "东方多策略灵活配置混合C"
$wsQ3.Cells.Item(2, 4).Value = "0.26"
$wsQ3.Cells.Item(2, 5).Value = "55.14"
$wsQ3.Cells.Item(2, 6).Value = "2.54"
$wsQ3.Cells.Item(2, 7).Value = "0.0066"
$wsQ3.Cells.Item(2, 8).Value = 5

# Row 3
$wsQ3.Cells.Item(3, 1).Value = 1
$wsQ3.Cells.Item(3, 2).Value = "015641"
$wsQ3.Cells.Item(3, 3).Value = "银华数字经济股票A"
$wsQ3.Cells.Item(3, 4).Value = "0.14"
$wsQ3.Cells.Item(3, 5).Value = "90.09"
$wsQ3.Cells.Item(3, 6).Value = "2.87"
$wsQ3.Cells.Item(3, 7).Value = "0.0040"
$wsQ3.Cells.Item(3, 8).Value = 7

# Row 4
$wsQ3.Cells.Item(4, 1).Value = 2
$wsQ3.Cells.Item(4, 2).Value = "400023"
$wsQ3.Cells.Item(4, 3).Value = "东方多策略灵活配置混合A"
$wsQ3.Cells.Item(4, 4).Value = "0.03"
$wsQ3.Cells.Item(4, 5).Value = "55.14"
$wsQ3.Cells.Item(4, 6).Value = "2.54"
$wsQ3.Cells.Item(4, 7).Value = "0.0008"
$wsQ3.Cells.Item(4, 8).Value = 5

# Row 5
$wsQ3.Cells.Item(5, 1).Value = 3
$wsQ3.Cells.Item(5, 2).Value = "015642"
$wsQ3.Cells.Item(5, 3).Value = "银华数字经济股票C"
$wsQ3.Cells.Item(5, 4).Value = "0.01"
$wsQ3.Cells.Item(5, 5).Value = "90.09"
$wsQ3.Cells.Item(5, 6).Value = "2.87"
$wsQ3.Cells.Item(5, 7).Value = "0.0003"
$wsQ3.Cells.Item(5, 8).Value = 7

# ---------------------------------------------------------------------------
# 2) Update the "总计" sheet: existing row 2 now refers to 2022-Q3, and a new
#    row 3 is appended with the original 2022-Q2 summary values.
# ---------------------------------------------------------------------------
$wsTotal.Cells.Item(2, 2).Value = "2022-Q3"

$wsTotal.Range("A2").Copy()
$wsTotal.Range("A3").PasteSpecial(-4122)

$wsTotal.Cells.Item(3, 1).Value = 1
$wsTotal.Cells.Item(3, 2).Value = "2022-Q2"
$wsTotal.Cells.Item(3, 3).Value = 4
$wsTotal.Cells.Item(3, 4).Value = 0.01

# ---------------------------------------------------------------------------
# 3) Restore the original active-tab selection (the "2022-Q2" sheet was the
#    selected tab before this edit); re-fetch it fresh by name since its
#    positional index changed when the new sheet was inserted before it.
# ---------------------------------------------------------------------------
$wsQ2Fresh = $wb.Worksheets.Item("2022-Q2")
$wsQ2Fresh.Activate()
